# Applies the scheduled Sheets update: refreshed market-price snapshot values
# across the per-job Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each entry below is one changed data row: sheet name, row number, and the
# column letter -> new numeric value map (current/average prices, leve costs,
# and the derived profit figures that depend on them). A $null value means the
# cell is cleared entirely (no longer populated for that row).

$wb = $excel.ActiveWorkbook

$updates = @(
  @{Sheet="ALC"; Row=15; Cells=@{"H"=972.2308; "I"=972.2308; "K"=2916.6924; "M"=-2747.6924}},
  @{Sheet="ALC"; Row=17; Cells=@{"H"=401735.2; "J"=401735.2; "L"=1205205.6; "N"=-1205541.6}},
  @{Sheet="ALC"; Row=41; Cells=@{"H"=661.7778; "I"=815.5833; "J"=354.16666; "K"=815.5833; "L"=354.16666; "M"=-375.5833; "N"=-1234.16666}},
  @{Sheet="ALC"; Row=53; Cells=@{"H"=106.333336; "I"=83; "J"=135.5; "K"=83; "L"=135.5; "M"=554; "N"=-1409.5}},
  @{Sheet="ALC"; Row=70; Cells=@{"H"=3132.182; "J"=2497.5715; "L"=7492.7145; "N"=-8032.7145}},
  @{Sheet="ALC"; Row=73; Cells=@{"H"=3132.182; "J"=2497.5715; "L"=7492.7145; "N"=-9364.7145}},
  @{Sheet="ALC"; Row=92; Cells=@{"H"=437.9643; "I"=382.9524; "J"=603; "K"=382.9524; "L"=603; "M"=865.0476; "N"=-3099}},
  @{Sheet="ALC"; Row=94; Cells=@{"H"=1370.375; "I"=1370.375; "K"=1370.375; "M"=-919.375}},
  @{Sheet="ALC"; Row=96; Cells=@{"H"=780.5294; "I"=738.13336; "K"=2214.40008; "M"=-841.4000800000003}},
  @{Sheet="ALC"; Row=98; Cells=@{"H"=2183.5625; "I"=1067.0714; "J"=9999; "K"=1067.0714; "L"=9999; "M"=430.9286; "N"=-12995}},
  @{Sheet="ALC"; Row=112; Cells=@{"H"=1496.2632; "J"=1520.5883; "L"=4561.7649; "N"=-6777.7649}},
  @{Sheet="ALC"; Row=122; Cells=@{"H"=2183.5625; "I"=1067.0714; "J"=9999; "K"=3201.2142; "L"=29997; "M"=-751.2142000000003; "N"=-34897}},
  @{Sheet="ALC"; Row=127; Cells=@{"H"=943.8333; "I"=932.6; "K"=2797.8; "M"=2162.2}},
  @{Sheet="ALC"; Row=129; Cells=@{"H"=1954.1; "I"=817.625; "J"=6500; "K"=2452.875; "L"=19500; "M"=2547.125; "N"=-29500}},
  @{Sheet="ALC"; Row=132; Cells=@{"H"=1807.5135; "I"=1834.6765; "K"=5504.029500000001; "M"=-2974.029500000001}},
  @{Sheet="ALC"; Row=138; Cells=@{"H"=2262.879; "I"=1854.1538; "J"=2528.55; "K"=5562.4614; "L"=7585.650000000001; "M"=-422.4614000000001; "N"=-17865.65}},
  @{Sheet="ALC"; Row=141; Cells=@{"H"=2539.3333; "I"=2539.3333; "K"=7617.999899999999; "M"=-2437.999899999999}},
  @{Sheet="ARM"; Row=19; Cells=@{"H"=6149.75; "I"=299.5; "J"=12000; "K"=299.5; "L"=12000; "M"=-70.5; "N"=-12458}},
  @{Sheet="ARM"; Row=21; Cells=@{"H"=2082.6667; "I"=1124; "J"=4000; "K"=1124; "L"=4000; "M"=-750; "N"=-4748}},
  @{Sheet="ARM"; Row=27; Cells=@{"H"=2000; "J"=2000; "L"=2000; "N"=-2368}},
  @{Sheet="ARM"; Row=30; Cells=@{"H"=1518; "J"=1366.3334; "L"=1366.3334; "N"=-1666.3334}},
  @{Sheet="ARM"; Row=32; Cells=@{"H"=9411.137000000001; "I"=3636.3416; "J"=88333.336; "K"=3636.3416; "L"=88333.336; "M"=-3349.3416; "N"=-88907.336}},
  @{Sheet="ARM"; Row=57; Cells=@{"H"=5136.25; "I"=5136.25; "K"=5136.25; "M"=-4652.25}},
  @{Sheet="ARM"; Row=63; Cells=@{"H"=1229.125; "I"=1229.125; "K"=1229.125; "M"=-543.125}},
  @{Sheet="ARM"; Row=66; Cells=@{"H"=1229.125; "I"=1229.125; "K"=6145.625; "M"=-2713.625}},
  @{Sheet="ARM"; Row=122; Cells=@{"H"=2166.182; "I"=2011.5; "K"=6034.5; "M"=-3584.5}},
  @{Sheet="BSM"; Row=20; Cells=@{"H"=7718.643; "I"=7673.9033; "J"=7844.727; "K"=7673.9033; "L"=7844.727; "M"=-7426.9033; "N"=-8338.726999999999}},
  @{Sheet="BSM"; Row=70; Cells=@{"H"=0; "J"=0; "L"=0; "N"=$null}},
  @{Sheet="BSM"; Row=73; Cells=@{"H"=0; "J"=0; "L"=0; "N"=$null}},
  @{Sheet="BSM"; Row=86; Cells=@{"H"=2716; "I"=2599.6; "K"=2599.6; "M"=-1476.6}},
  @{Sheet="BSM"; Row=89; Cells=@{"H"=2716; "I"=2599.6; "K"=12998; "M"=-7382}},
  @{Sheet="BSM"; Row=107; Cells=@{"H"=85094.5; "I"=112492.664; "K"=112492.664; "M"=-110572.664}},
  @{Sheet="BSM"; Row=113; Cells=@{"H"=5092.5454; "I"=5092.5454; "K"=5092.5454; "M"=-2922.5454}},
  @{Sheet="BSM"; Row=128; Cells=@{"H"=19626.834; "I"=19626.834; "K"=58880.50199999999; "M"=-56390.50199999999}},
  @{Sheet="BSM"; Row=134; Cells=@{"H"=1337.3611; "I"=1157.4814; "J"=1877; "K"=3472.4442; "L"=5631; "M"=-937.4441999999999; "N"=-10701}},
  @{Sheet="CRP"; Row=31; Cells=@{"H"=12632.243; "I"=3552.1667; "J"=29395.46; "K"=3552.1667; "L"=29395.46; "M"=-3257.1667; "N"=-29985.46}},
  @{Sheet="CRP"; Row=34; Cells=@{"H"=12632.243; "I"=3552.1667; "J"=29395.46; "K"=3552.1667; "L"=29395.46; "M"=-3350.1667; "N"=-29799.46}},
  @{Sheet="CRP"; Row=62; Cells=@{"H"=3100.3; "I"=3187.875; "J"=2750; "K"=3187.875; "L"=2750; "M"=-2563.875; "N"=-3998}},
  @{Sheet="CRP"; Row=65; Cells=@{"H"=3100.3; "I"=3187.875; "J"=2750; "K"=15939.375; "L"=13750; "M"=-12819.375; "N"=-19990}},
  @{Sheet="CRP"; Row=70; Cells=@{"H"=48567.5; "J"=48567.5; "L"=48567.5; "N"=-49197.5}},
  @{Sheet="CRP"; Row=73; Cells=@{"H"=48567.5; "J"=48567.5; "L"=48567.5; "N"=-50751.5}},
  @{Sheet="CRP"; Row=76; Cells=@{"H"=5070; "I"=5070; "K"=5070; "M"=-4755}},
  @{Sheet="CRP"; Row=79; Cells=@{"H"=5070; "I"=5070; "K"=5070; "M"=-3978}},
  @{Sheet="CRP"; Row=122; Cells=@{"H"=70311.87; "I"=112064.664; "K"=336193.992; "M"=-333743.992}},
  @{Sheet="CRP"; Row=132; Cells=@{"H"=3530.0667; "I"=3500.9614; "J"=3719.25; "K"=10502.8842; "L"=11157.75; "M"=-7972.8842; "N"=-16217.75}},
  @{Sheet="CRP"; Row=134; Cells=@{"H"=2477.5334; "I"=1553.4546; "J"=5018.75; "K"=4660.3638; "L"=15056.25; "M"=-2125.3638; "N"=-20126.25}},
  @{Sheet="CUL"; Row=56; Cells=@{"H"=9114.315000000001; "I"=9114.315000000001; "K"=9114.315000000001; "M"=-8584.315000000001}},
  @{Sheet="CUL"; Row=68; Cells=@{"H"=863; "J"=979; "L"=2937; "N"=-4559}},
  @{Sheet="CUL"; Row=71; Cells=@{"H"=863; "J"=979; "L"=8811; "N"=-16923}},
  @{Sheet="CUL"; Row=86; Cells=@{"H"=655.2174; "I"=467.27274; "K"=1401.81822; "M"=-215.8182200000001}},
  @{Sheet="CUL"; Row=89; Cells=@{"H"=655.2174; "I"=467.27274; "K"=4205.45466; "M"=1722.54534}},
  @{Sheet="CUL"; Row=98; Cells=@{"H"=535.5; "I"=688; "J"=281.33334; "K"=2064; "L"=844.0000200000001; "M"=-566; "N"=-3840.00002}},
  @{Sheet="CUL"; Row=132; Cells=@{"H"=2635.6365; "I"=2000; "J"=2776.889; "K"=18000; "L"=24992.001; "M"=-15470; "N"=-30052.001}},
  @{Sheet="CUL"; Row=134; Cells=@{"H"=13174.286; "I"=1934.2222; "J"=33406.4; "K"=5802.6666; "L"=100219.2; "M"=-732.6665999999996; "N"=-110359.2}},
  @{Sheet="GSM"; Row=11; Cells=@{"H"=2435811.5; "I"=1584702.9; "J"=3905908.2; "K"=1584702.9; "L"=3905908.2; "M"=-1584563.9; "N"=-3906186.2}},
  @{Sheet="GSM"; Row=18; Cells=@{"H"=0; "I"=0; "K"=0; "M"=$null}},
  @{Sheet="GSM"; Row=22; Cells=@{"H"=386.75; "I"=564; "J"=209.5; "K"=564; "L"=209.5; "M"=-35; "N"=-1267.5}},
  @{Sheet="GSM"; Row=27; Cells=@{"H"=1000; "J"=1000; "L"=1000; "N"=-1332}},
  @{Sheet="GSM"; Row=43; Cells=@{"H"=22091.4; "I"=3819.3333; "K"=3819.3333; "M"=-3668.3333}},
  @{Sheet="GSM"; Row=46; Cells=@{"H"=36999.766; "J"=35777.332; "L"=35777.332; "N"=-36089.332}},
  @{Sheet="GSM"; Row=80; Cells=@{"H"=4527.8; "I"=3216.3333; "K"=3216.3333; "M"=-2218.3333}},
  @{Sheet="GSM"; Row=83; Cells=@{"H"=4527.8; "I"=3216.3333; "K"=16081.6665; "M"=-11089.6665}},
  @{Sheet="GSM"; Row=122; Cells=@{"H"=3307.2856; "I"=2776.75; "K"=8330.25; "M"=-5880.25}},
  @{Sheet="GSM"; Row=126; Cells=@{"H"=2397; "I"=0; "J"=2397; "K"=0; "L"=7191; "M"=$null; "N"=-12131}},
  @{Sheet="GSM"; Row=132; Cells=@{"H"=5974.4443; "I"=5974.4443; "K"=17923.3329; "M"=-15393.3329}},
  @{Sheet="LTW"; Row=22; Cells=@{"H"=2000; "J"=2000; "L"=2000; "N"=-2590}},
  @{Sheet="LTW"; Row=27; Cells=@{"H"=2000; "J"=2000; "L"=2000; "N"=-2214}},
  @{Sheet="LTW"; Row=36; Cells=@{"H"=10715; "J"=10715; "L"=10715; "N"=-11839}},
  @{Sheet="LTW"; Row=40; Cells=@{"H"=10395.546; "I"=10721.158; "J"=8333.333000000001; "K"=10721.158; "L"=8333.333000000001; "M"=-10585.158; "N"=-8605.333000000001}},
  @{Sheet="LTW"; Row=62; Cells=@{"H"=0; "J"=0; "L"=0; "N"=$null}},
  @{Sheet="LTW"; Row=65; Cells=@{"H"=0; "J"=0; "L"=0; "N"=$null}},
  @{Sheet="LTW"; Row=82; Cells=@{"H"=1463.25; "I"=1531.4286; "J"=1410.2222; "K"=1531.4286; "L"=1410.2222; "M"=-1170.4286; "N"=-2132.2222}},
  @{Sheet="LTW"; Row=85; Cells=@{"H"=1463.25; "I"=1531.4286; "J"=1410.2222; "K"=1531.4286; "L"=1410.2222; "M"=-283.4286; "N"=-3906.2222}},
  @{Sheet="LTW"; Row=124; Cells=@{"H"=59429; "J"=59429; "L"=59429; "N"=-69249}},
  @{Sheet="LTW"; Row=132; Cells=@{"H"=3214.3225; "I"=2946.16; "J"=4331.6665; "K"=8838.48; "L"=12994.9995; "M"=-6308.48; "N"=-18054.9995}},
  @{Sheet="LTW"; Row=136; Cells=@{"H"=3921.9312; "I"=3635.9; "J"=4557.5557; "K"=10907.7; "L"=13672.6671; "M"=-8357.700000000001; "N"=-18772.6671}},
  @{Sheet="WVR"; Row=4; Cells=@{"H"=34999.5; "I"=40000; "K"=40000; "M"=-39887}},
  @{Sheet="WVR"; Row=32; Cells=@{"H"=0; "I"=0; "K"=0; "M"=$null}},
  @{Sheet="WVR"; Row=62; Cells=@{"H"=18573.455; "I"=3125.5; "K"=3125.5; "M"=-2501.5}},
  @{Sheet="WVR"; Row=65; Cells=@{"H"=18573.455; "I"=3125.5; "K"=15627.5; "M"=-12507.5}},
  @{Sheet="WVR"; Row=96; Cells=@{"H"=2117.9092; "I"=2229.7; "K"=2229.7; "M"=-856.6999999999998}},
  @{Sheet="WVR"; Row=119; Cells=@{"H"=27500; "J"=27500; "L"=27500; "N"=-37176}},
  @{Sheet="WVR"; Row=126; Cells=@{"H"=1861.1904; "I"=1849.75; "J"=1897.8; "K"=5549.25; "L"=5693.4; "M"=-3079.25; "N"=-10633.4}},
  @{Sheet="WVR"; Row=132; Cells=@{"H"=1524.4722; "I"=1002.913; "K"=3008.739; "M"=-478.739}}
)

foreach ($update in $updates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$($update.Row)"
        $value = $update.Cells[$col]
        if ($null -eq $value) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $value
        }
    }
}
